# Update column G ("K") values for rows 2-10 on the active worksheet.
# This mirrors the save_data regeneration described in the commit message:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
